$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 70: new Coding time-log entry
$ws.Cells.Item(70, 1).Value = 41913
$ws.Cells.Item(70, 2).Value = 0.52777777777777779
$ws.Cells.Item(70, 3).Value = 0.57777777777777783
$ws.Cells.Item(70, 4).Value = 15
$ws.Cells.Item(70, 6).Value = "Coding"
$ws.Cells.Item(70, 5).Formula = "=IF(AND(NOT(ISBLANK(B70)),NOT(ISBLANK(C70))),(C70-B70)*24-D70/60,"""")"

# Row 71: new Coding time-log entry
$ws.Cells.Item(71, 1).Value = 41913
$ws.Cells.Item(71, 2).Value = 0.8881944444444444
$ws.Cells.Item(71, 3).Value = 1.0756944444444445
$ws.Cells.Item(71, 4).Value = 30
$ws.Cells.Item(71, 6).Value = "Coding"
$ws.Cells.Item(71, 5).Formula = "=IF(AND(NOT(ISBLANK(B71)),NOT(ISBLANK(C71))),(C71-B71)*24-D71/60,"""")"

# Move the active selection, matching where the user ended up working next
[void]$ws.Range("D72").Select()
